$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.410.55'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.567.82'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''1.000'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '''286.84'
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = '''0.3749'
$ws.Range("E7").Value = '  +2.90%  '
$ws.Range("D8").Value = '''0.3278'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").Value = '''45.53'
$ws.Range("E9").Value = '  -5.18%  '
$ws.Range("D10").Value = '''1.151'
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("D11").Value = '''0.07429'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '''1.000'
$ws.Range("D13").Value = '''20.51'
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").Value = '''5.852'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").Value = '''6.841'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '1.571.64'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").Value = '''0.00001101'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '''0.06698'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '''86.04'
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = '''0.9999'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '''6.362'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("D23").Value = '''11.71'
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").Value = '22.399.92'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '''2.319'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = '''2.570'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '''152.04'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").Value = '''19.34'
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = '''4.919'
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").Value = '''123.51'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '1.747.88'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").Value = '''1.062'
$ws.Range("E32").Value = '  +3.80%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.946'
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '''1.946'
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("D35").Value = '''9.672'
$ws.Range("E35").Value = '  -0.68%  '
$ws.Range("D36").Value = '''0.08265'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '''0.02389'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = '''1.291'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '''0.06332'
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").Value = '''0.2195'
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("D41").Value = '''5.284'
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '''0.6115'
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '''0.9997'
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''13.70'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").Value = '''3.750'
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5926'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''2.013'
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''124.14'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '''1.182'
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.07156'
$ws.Range("E51").Value = '  -0.61%  '
